$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. "326.04") are not
# auto-converted to numbers by the Value setter, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.194.76"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.839.56"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "326.04"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "0.4639"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "0.3870"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "0.07865"
$ws.Range("D10").Value = "0.9638"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.864.66"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "6.880"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "0.06875"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "88.54"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "0.000009973"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "28.202.31"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "5.304"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "2.102"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "2.057.49"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "154.74"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "19.15"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "5.727"
$ws.Range("E28").Value = "  -5.43%  "
$ws.Range("D29").Value = "1.971"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "119.03"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "0.09261"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "0.9322"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").Value = "5.289"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "3.332"
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("D36").Value = "0.05829"
$ws.Range("E36").Value = "  -4.54%  "
$ws.Range("D37").Value = "0.02128"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").Value = "1.138"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "7.765"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "0.5599"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "9.896"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").Value = "0.1762"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "0.07273"
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "11.62"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "0.5278"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "1.137"
$ws.Range("E46").Value = "  -8.73%  "
$ws.Range("D47").Value = "2.133"
$ws.Range("E47").Value = "  -11.94%  "
$ws.Range("D48").Value = "1.836"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").Value = "114.04"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "2.326"
$ws.Range("E51").Value = "  -0.74%  "

# Restore the original (default/general) formatting so cell styles are
# unchanged from the source workbook.
$ws.Range("D2:E51").NumberFormat = "General"
$ws.Range("D2:E51").Style = "Normal"
